{"js": "// Update the worksheet date heading and all 100 arithmetic-problem cells\n// in the single table, per the target revision.\n\nconst NEW_DATE = \"2024-10-15 Tuesday\";\nconst OLD_DATE = \"2024-10-14 Monday\";\n\n// New contents for the 20 x 5 table, row-major (top-left to bottom-right).\nconst NEW_GRID = [\n  [\"71+8=\", \"92-20=\", \"42-20=\", \"84-49=\", \"53+2=\"],\n  [\"56-21=\", \"74-61=\", \"49-8=\", \"76-74=\", \"2+18=\"],\n  [\"93-45=\", \"2+29=\", \"14+55=\", \"89-23=\", \"42+6=\"],\n  [\"63+7=\", \"86-9=\", \"25+28=\", \"29-22=\", \"76-62=\"],\n  [\"40+7=\", \"19+34=\", \"90-54=\", \"66+14=\", \"69-64=\"],\n  [\"39+29=\", \"12-2=\", \"39+39=\", \"46+5=\", \"86-5=\"],\n  [\"29-5=\", \"77-37=\", \"63+27=\", \"78-61=\", \"60-2=\"],\n  [\"12+1=\", \"86-45=\", \"54+36=\", \"85-4=\", \"97-27=\"],\n  [\"35+28=\", \"6+4=\", \"30+9=\", \"82-44=\", \"37+57=\"],\n  [\"73-31=\", \"35+20=\", \"39-38=\", \"46-5=\", \"36-4=\"],\n  [\"13+69=\", \"64+4=\", \"10-9=\", \"3+87=\", \"68-47=\"],\n  [\"48-25=\", \"64-35=\", \"50-7=\", \"83-43=\", \"14+47=\"],\n  [\"29-18=\", \"43+28=\", \"44+53=\", \"85-42=\", \"74+23=\"],\n  [\"74-42=\", \"19+41=\", \"53+33=\", \"74+24=\", \"3+86=\"],\n  [\"18+1=\", \"49-36=\", \"73+22=\", \"92-30=\", \"64-56=\"],\n  [\"90-14=\", \"3+95=\", \"84-31=\", \"50-3=\", \"89-79=\"],\n  [\"91-43=\", \"52-0=\", \"55-19=\", \"22+9=\", \"75-29=\"],\n  [\"7+31=\", \"83-2=\", \"59+13=\", \"49-26=\", \"63-50=\"],\n  [\"55+38=\", \"64-45=\", \"66-1=\", \"8+89=\", \"33+27=\"],\n  [\"24-11=\", \"68+14=\", \"40+41=\", \"83-9=\", \"59-1=\"],\n];\n\n// 1) Update the date paragraph above the table, preserving its run formatting.\nconst dateResults = context.document.body.search(OLD_DATE, { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(NEW_DATE, Word.InsertLocation.replace);\n} else {\n  // Fallback: if the exact old date text isn't found (e.g. already edited),\n  // update the first paragraph directly.\n  const paragraphs = context.document.body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n  if (paragraphs.items.length > 0) {\n    paragraphs.items[0].insertText(NEW_DATE, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n\n// 2) Update every cell of the (single) table with the new problem text,\n//    using Cell.value so existing run/paragraph formatting is kept.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nfor (let r = 0; r < table.rowCount; r++) {\n  const rowValues = NEW_GRID[r];\n  if (!rowValues) continue;\n  for (let c = 0; c < rowValues.length; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = rowValues[c];\n  }\n}\nawait context.sync();\n", "ps1": "# Update the worksheet date heading and all 100 arithmetic-problem cells\n# in the single table, per the target revision.\n\n$d = $word.ActiveDocument\n\n$oldDate = \"2024-10-14 Monday\"\n$newDate = \"2024-10-15 Tuesday\"\n\n# 1) Update the date paragraph above the table, preserving its run formatting.\n$find = $d.Content.Find\n$find.Text = $oldDate\n$find.Replacement.Text = $newDate\n$find.Execute($oldDate, $false, $false, $false, $false, $false, $true, 1, $false, $newDate, 2)\n\n# New contents for the 20 x 5 table, row-major (top-left to bottom-right).\n$newGrid = @(\n    @(\"71+8=\", \"92-20=\", \"42-20=\", \"84-49=\", \"53+2=\"),\n    @(\"56-21=\", \"74-61=\", \"49-8=\", \"76-74=\", \"2+18=\"),\n    @(\"93-45=\", \"2+29=\", \"14+55=\", \"89-23=\", \"42+6=\"),\n    @(\"63+7=\", \"86-9=\", \"25+28=\", \"29-22=\", \"76-62=\"),\n    @(\"40+7=\", \"19+34=\", \"90-54=\", \"66+14=\", \"69-64=\"),\n    @(\"39+29=\", \"12-2=\", \"39+39=\", \"46+5=\", \"86-5=\"),\n    @(\"29-5=\", \"77-37=\", \"63+27=\", \"78-61=\", \"60-2=\"),\n    @(\"12+1=\", \"86-45=\", \"54+36=\", \"85-4=\", \"97-27=\"),\n    @(\"35+28=\", \"6+4=\", \"30+9=\", \"82-44=\", \"37+57=\"),\n    @(\"73-31=\", \"35+20=\", \"39-38=\", \"46-5=\", \"36-4=\"),\n    @(\"13+69=\", \"64+4=\", \"10-9=\", \"3+87=\", \"68-47=\"),\n    @(\"48-25=\", \"64-35=\", \"50-7=\", \"83-43=\", \"14+47=\"),\n    @(\"29-18=\", \"43+28=\", \"44+53=\", \"85-42=\", \"74+23=\"),\n    @(\"74-42=\", \"19+41=\", \"53+33=\", \"74+24=\", \"3+86=\"),\n    @(\"18+1=\", \"49-36=\", \"73+22=\", \"92-30=\", \"64-56=\"),\n    @(\"90-14=\", \"3+95=\", \"84-31=\", \"50-3=\", \"89-79=\"),\n    @(\"91-43=\", \"52-0=\", \"55-19=\", \"22+9=\", \"75-29=\"),\n    @(\"7+31=\", \"83-2=\", \"59+13=\", \"49-26=\", \"63-50=\"),\n    @(\"55+38=\", \"64-45=\", \"66-1=\", \"8+89=\", \"33+27=\"),\n    @(\"24-11=\", \"68+14=\", \"40+41=\", \"83-9=\", \"59-1=\")\n)\n\n# 2) Update every cell of the (single) table with the new problem text,\n#    using Cell.Range.Text so existing run/paragraph formatting is kept.\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n    $rowValues = $newGrid[$r - 1]\n    for ($c = 1; $c -le $colCount; $c++) {\n        $t.Cell($r, $c).Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
